$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 65 (Feria Lagunitas de Puerto Montt,
# Albahaca data set), shifting the existing rows 65-144 down to 66-145.
$ws.Rows(65).Insert()

# Populate the newly inserted row 65 with the new weekly price record.
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44848
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = 100112052
$ws.Range("G65").Value = "Albahaca"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 90
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = 7500
$ws.Range("N65").Value = "$/paquete"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 7500
$ws.Range("Q65").Value = 1
$ws.Range("R65").Value = "Hortaliza"
